$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.05070303331908121
$ws.Range("D2").Value = 0.2332309276366971
$ws.Range("E2").Value = 0.04604448509146364
$ws.Range("F2").Value = 14.91324311680154
$ws.Range("G2").Value = 0.002918073986712848
$ws.Range("I2").Value = 11.19940568461726
$ws.Range("J2").Value = 0.4134903892802129
$ws.Range("L2").Value = 0.2071373558207199
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("C3").Value = 0.04543810182417474
$ws.Range("D3").Value = 0.2199300229874552
$ws.Range("E3").Value = 0.0462045057043019
$ws.Range("F3").Value = 14.716383197151
$ws.Range("G3").Value = 0.00293686367382925
$ws.Range("I3").Value = 11.04696120098782
$ws.Range("J3").Value = 0.4132330119170931
$ws.Range("L3").Value = 0.2084080281729008
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("C4").Value = 0.04225366310134859
$ws.Range("D4").Value = 0.2119966071514625
$ws.Range("E4").Value = 0.04630899095269991
$ws.Range("F4").Value = 14.60697997074351
$ws.Range("G4").Value = 0.002948941771239945
$ws.Range("I4").Value = 10.96202227129362
$ws.Range("J4").Value = 0.4133321067961759
$ws.Range("L4").Value = 0.2092837934835892
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("C5").Value = 0.04096734721807138
$ws.Range("D5").Value = 0.2088201796358078
$ws.Range("E5").Value = 0.04635313975006827
$ws.Range("F5").Value = 14.56523231834893
$ws.Range("G5").Value = 0.002954000779421528
$ws.Range("I5").Value = 10.92954969324299
$ws.Range("J5").Value = 0.4134365316262603
$ws.Range("L5").Value = 0.2096646463383678
$ws.Range("N5").Value = 2.293303068605894
$ws.Range("C6").Value = 0.04075442023130904
$ws.Range("D6").Value = 0.2082960865213863
$ws.Range("E6").Value = 0.04636056555516221
$ws.Range("F6").Value = 14.55846998389433
$ws.Range("G6").Value = 0.002954849132652956
$ws.Range("I6").Value = 10.92428587027445
$ws.Range("J6").Value = 0.413457721899178
$ws.Range("L6").Value = 0.2097293328199683
$ws.Range("N6").Value = 2.2803432614038
$ws.Range("C7").Value = 0.04223627035509026
$ws.Range("D7").Value = 0.2119535427891606
$ws.Range("E7").Value = 0.04630957999684027
$ws.Range("F7").Value = 14.60640553169947
$ws.Range("G7").Value = 0.002949009442290111
$ws.Range("I7").Value = 10.96157571557316
$ws.Range("J7").Value = 0.413333256588146
$ws.Range("L7").Value = 0.2092888328137334
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("C8").Value = 0.04887720126328077
$ws.Range("D8").Value = 0.2285949924482793
$ws.Range("E8").Value = 0.04609836862104189
$ws.Range("F8").Value = 14.84295510381548
$ws.Range("G8").Value = 0.002924440968055497
$ws.Range("I8").Value = 11.14502178447376
$ws.Range("J8").Value = 0.4133479004539282
$ws.Range("L8").Value = 0.2075556110779466
$ws.Range("N8").Value = 2.766433886209882
$ws.Range("C9").Value = 0.06231916246990465
$ws.Range("D9").Value = 0.263184562747341
$ws.Range("E9").Value = 0.04573349430175133
$ws.Range("F9").Value = 15.40017927948475
$ws.Range("G9").Value = 0.002880509321703251
$ws.Range("I9").Value = 11.57528781273442
$ws.Range("J9").Value = 0.4154463219016691
$ws.Range("L9").Value = 0.2049180654707214
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("C10").Value = 0.07250352962753936
$ws.Range("D10").Value = 0.2899419851918594
$ws.Range("E10").Value = 0.04549529974901456
$ws.Range("F10").Value = 15.86993682970268
$ws.Range("G10").Value = 0.002850755127416633
$ws.Range("I10").Value = 11.9370448425725
$ws.Range("J10").Value = 0.4182935905188856
$ws.Range("L10").Value = 0.2034490359930601
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("C11").Value = 0.07721590248792154
$ws.Range("D11").Value = 0.3024417553300225
$ws.Range("E11").Value = 0.04539339087058236
$ws.Range("F11").Value = 16.09755312645564
$ws.Range("G11").Value = 0.002837751999561609
$ws.Range("I11").Value = 12.11214286542514
$ws.Range("J11").Value = 0.419882658873604
$ws.Range("L11").Value = 0.2028836782583454
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("C12").Value = 0.07901282628256467
$ws.Range("D12").Value = 0.30722522566316
$ws.Range("E12").Value = 0.0453557252531045
$ws.Range("F12").Value = 16.18581648607648
$ws.Range("G12").Value = 0.002832903357964923
$ws.Range("I12").Value = 12.18001544486083
$ws.Range("J12").Value = 0.4205275429961972
$ws.Range("L12").Value = 0.2026844950552231
$ws.Range("N12").Value = 4.460285735713398
$ws.Range("C13").Value = 0.07862525599406922
$ws.Range("D13").Value = 0.3061927437985901
$ws.Range("E13").Value = 0.04536379610535901
$ws.Range("F13").Value = 16.16671420453912
$ws.Range("G13").Value = 0.002833944266697669
$ws.Range("I13").Value = 12.16532732895178
$ws.Range("J13").Value = 0.4203867224228901
$ws.Range("L13").Value = 0.2027267280018776
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("C14").Value = 0.07736348091248146
$ws.Range("D14").Value = 0.3028342731927296
$ws.Range("E14").Value = 0.04539027356878123
$ws.Range("F14").Value = 16.10477271695669
$ws.Range("G14").Value = 0.002837351595314271
$ws.Range("I14").Value = 12.11769507480591
$ws.Range("J14").Value = 0.4199348440663329
$ws.Range("L14").Value = 0.2028669919647186
$ws.Range("N14").Value = 4.371107314139238
$ws.Range("C15").Value = 0.07659225993793939
$ws.Range("D15").Value = 0.3007837249947158
$ws.Range("E15").Value = 0.0454066121936938
$ws.Range("F15").Value = 16.06710342733282
$ws.Range("G15").Value = 0.002839448460991947
$ws.Range("I15").Value = 12.08872457605707
$ws.Range("J15").Value = 0.4196637005733095
$ws.Range("L15").Value = 0.2029548521475704
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("C16").Value = 0.07219725028929247
$ws.Range("D16").Value = 0.2891319392111882
$ws.Range("E16").Value = 0.04550208928837152
$ws.Range("F16").Value = 15.85534744512699
$ws.Range("G16").Value = 0.002851615522453562
$ws.Range("I16").Value = 11.92581811945206
$ws.Range("J16").Value = 0.4181957283879427
$ws.Range("L16").Value = 0.203488063243455
$ws.Range("N16").Value = 4.089429168003846
$ws.Range("C17").Value = 0.06952211806698472
$ws.Range("D17").Value = 0.2820699329581657
$ws.Range("E17").Value = 0.04556231114008913
$ws.Range("F17").Value = 15.72905494878484
$ws.Range("G17").Value = 0.002859215109046806
$ws.Range("I17").Value = 11.82861426732723
$ws.Range("J17").Value = 0.4173710468893148
$ws.Range("L17").Value = 0.2038416025143022
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("C18").Value = 0.06799086703441048
$ws.Range("D18").Value = 0.2780387114119094
$ws.Range("E18").Value = 0.04559755608920568
$ws.Range("F18").Value = 15.65772089075125
$ws.Range("G18").Value = 0.002863636355994122
$ws.Range("I18").Value = 11.77369331330868
$ws.Range("J18").Value = 0.4169243127736308
$ws.Range("L18").Value = 0.2040546285414635
$ws.Range("N18").Value = 3.857331695637072
$ws.Range("C19").Value = 0.0674736547126713
$ws.Range("D19").Value = 0.276678986775579
$ws.Range("E19").Value = 0.04560959372901419
$ws.Range("F19").Value = 15.63379069995256
$ws.Range("G19").Value = 0.002865141965192757
$ws.Range("I19").Value = 11.75526618979336
$ws.Range("J19").Value = 0.4167777705747682
$ws.Range("L19").Value = 0.2041284148472258
$ws.Range("N19").Value = 3.828614786363971
$ws.Range("C20").Value = 0.06980611584307894
$ws.Range("D20").Value = 0.2828185016638827
$ws.Range("E20").Value = 0.04555583762105497
$ws.Range("F20").Value = 15.74236334751231
$ws.Range("G20").Value = 0.002858400936445182
$ws.Range("I20").Value = 11.83885914654104
$ws.Range("J20").Value = 0.4174559728007381
$ws.Range("L20").Value = 0.2038029651063553
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("C21").Value = 0.07773374851836934
$ws.Range("D21").Value = 0.3038193539928784
$ws.Range("E21").Value = 0.04538247140498974
$ws.Range("F21").Value = 16.12290969293531
$ws.Range("G21").Value = 0.002836348744423452
$ws.Range("I21").Value = 12.1316428719735
$ws.Range("J21").Value = 0.4200663936182423
$ws.Range("L21").Value = 0.2028253876126414
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("C22").Value = 0.08298799772418874
$ws.Range("D22").Value = 0.3178380017768347
$ws.Range("E22").Value = 0.04527455774199352
$ws.Range("F22").Value = 16.38371855949526
$ws.Range("G22").Value = 0.002822375055057358
$ws.Range("I22").Value = 12.33215241920033
$ws.Range("J22").Value = 0.4220243461837185
$ws.Range("L22").Value = 0.202273420662948
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("C23").Value = 0.0801766670385291
$ws.Range("D23").Value = 0.3103281254266221
$ws.Range("E23").Value = 0.045331660615612
$ws.Range("F23").Value = 16.24338935905934
$ws.Range("G23").Value = 0.002829793331219489
$ws.Range("I23").Value = 12.22428075172343
$ws.Range("J23").Value = 0.4209559959261924
$ws.Range("L23").Value = 0.2025600235459919
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("C24").Value = 0.06967769959950942
$ws.Range("D24").Value = 0.2824799842062475
$ws.Range("E24").Value = 0.04555876235954059
$ws.Range("F24").Value = 15.73634265988687
$ws.Range("G24").Value = 0.002858768861464047
$ws.Range("I24").Value = 11.83422444154195
$ws.Range("J24").Value = 0.417417492560304
$ws.Range("L24").Value = 0.2038204026507557
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("C25").Value = 0.05863254563145404
$ws.Range("D25").Value = 0.2536024829851726
$ws.Range("E25").Value = 0.04582694284855893
$ws.Range("F25").Value = 15.2390946484964
$ws.Range("G25").Value = 0.002891946121026568
$ws.Range("I25").Value = 11.45107342679313
$ws.Range("J25").Value = 0.4146526848918199
$ws.Range("L25").Value = 0.2055496314698378
$ws.Range("N25").Value = 3.331249627311138
